# Refresh cryptos price/volume snapshot (scheduled GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.329.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "1.860.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'323.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  -6.25%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3865"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.14%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'48.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.62%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.07902"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.08%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'21.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "1.876.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.889"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'7.152"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.94%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001033"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'85.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.92%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06524"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'17.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.517"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.85%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "27.338.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "  -6.09%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.266"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "2.084.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'152.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'19.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "  -5.68%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'5.512"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.90%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'120.34"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = "'1.493"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.09308"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.9371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.600"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'5.286"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.50%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.02233"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.06002"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.214"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'8.269"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.65%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.9996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.5909"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1890"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'10.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.71%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.5653"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.44%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'11.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.61%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.03%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'3.366"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.06798"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = "'108.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.01%  "
$ws.Range("E51").Style = "Normal"

